$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 26 (which currently holds "CORNER"),
# shifting CORNER and everything below it down by one row.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new CHTRANS command entry.
$ws.Range("A26").Value = "CHTRANS"
$ws.Range("B26").Value = "Changes all the entities in a selected block to a desired transparency"

# Match the author's new selection position.
$ws.Range("B27").Select()
